$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Price and Volume columns to Text format so numeric-looking strings
# (e.g. "31.479.50", "0.000007895", "25.80") are preserved exactly as text.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '31.479.50'
$ws.Range("E2").Value = '  +3.52%  '

$ws.Range("D3").Value = '1.989.60'
$ws.Range("E3").Value = '  +5.76%  '

$ws.Range("D4").Value = '0.9955'
$ws.Range("E4").Value = '  -0.68%  '

$ws.Range("D5").Value = '0.7948'
$ws.Range("E5").Value = '  +68.94%  '

$ws.Range("D6").Value = '253.09'
$ws.Range("E6").Value = '  +3.95%  '

$ws.Range("D7").Value = '0.9954'
$ws.Range("E7").Value = '  -0.57%  '

$ws.Range("D8").Value = '0.3435'
$ws.Range("E8").Value = '  +19.51%  '

$ws.Range("D9").Value = '25.80'
$ws.Range("E9").Value = '  +16.98%  '

$ws.Range("D10").Value = '0.06974'
$ws.Range("E10").Value = '  +8.30%  '

$ws.Range("D11").Value = '0.8452'
$ws.Range("E11").Value = '  +17.40%  '

$ws.Range("D12").Value = '0.08129'
$ws.Range("E12").Value = '  +4.58%  '

$ws.Range("D13").Value = '102.63'
$ws.Range("E13").Value = '  +7.81%  '

$ws.Range("D14").Value = '1.982.29'
$ws.Range("E14").Value = '  +5.39%  '

$ws.Range("D15").Value = '5.512'
$ws.Range("E15").Value = '  +6.92%  '

$ws.Range("D16").Value = '276.03'
$ws.Range("E16").Value = '  -1.00%  '

$ws.Range("D17").Value = '31.379.30'
$ws.Range("E17").Value = '  +3.14%  '

$ws.Range("D18").Value = '14.01'
$ws.Range("E18").Value = '  +7.99%  '

$ws.Range("D19").Value = '0.000007895'
$ws.Range("E19").Value = '  +6.25%  '

$ws.Range("D20").Value = '2.240.58'
$ws.Range("E20").Value = '  +4.97%  '

$ws.Range("D21").Value = '5.696'
$ws.Range("E21").Value = '  +8.67%  '

$ws.Range("D22").Value = '0.9988'
$ws.Range("E22").Value = '  -0.20%  '

$ws.Range("D23").Value = '0.9965'
$ws.Range("E23").Value = '  -0.70%  '

$ws.Range("D24").Value = '6.875'
$ws.Range("E24").Value = '  +9.81%  '

$ws.Range("B25").Value = 'Cosmos'
$ws.Range("C25").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D25").Value = '9.671'
$ws.Range("E25").Value = '  +7.11%  '

$ws.Range("B26").Value = 'Stellar'
$ws.Range("C26").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D26").Value = '0.1559'
$ws.Range("E26").Value = '  +61.94%  '

$ws.Range("D27").Value = '166.49'
$ws.Range("E27").Value = '  +1.61%  '

$ws.Range("D28").Value = '19.67'
$ws.Range("E28").Value = '  +4.92%  '

$ws.Range("D29").Value = '2.295'
$ws.Range("E29").Value = '  +22.07%  '

$ws.Range("D30").Value = '1.563'
$ws.Range("E30").Value = '  +6.72%  '

$ws.Range("D31").Value = '1.352'
$ws.Range("E31").Value = '  +1.28%  '

$ws.Range("D32").Value = '4.573'
$ws.Range("E32").Value = '  +7.86%  '

$ws.Range("D33").Value = '4.338'
$ws.Range("E33").Value = '  +5.42%  '

$ws.Range("D34").Value = '0.05213'
$ws.Range("E34").Value = '  +7.77%  '

$ws.Range("D35").Value = '1.219'
$ws.Range("E35").Value = '  +8.84%  '

$ws.Range("D36").Value = '0.7459'
$ws.Range("E36").Value = '  +8.30%  '

$ws.Range("D37").Value = '2.794'
$ws.Range("E37").Value = '  +3.02%  '

$ws.Range("B38").Value = 'Frax'
$ws.Range("C38").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D38").Value = '0.9960'
$ws.Range("E38").Value = '  -0.44%  '

$ws.Range("B39").Value = 'VeChain'
$ws.Range("C39").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D39").Value = '0.01985'
$ws.Range("E39").Value = '  +5.77%  '

$ws.Range("B40").Value = 'MXToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D40").Value = '2.906'
$ws.Range("E40").Value = '  +3.24%  '

$ws.Range("B41").Value = 'FraxShare'
$ws.Range("C41").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D41").Value = '6.619'
$ws.Range("E41").Value = '  +7.34%  '

$ws.Range("B42").Value = 'Aave'
$ws.Range("C42").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D42").Value = '78.62'
$ws.Range("E42").Value = '  +5.80%  '

$ws.Range("B43").Value = 'TheSandbox'
$ws.Range("C43").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D43").Value = '0.4678'
$ws.Range("E43").Value = '  +10.64%  '

$ws.Range("B44").Value = 'RenderToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D44").Value = '2.084'
$ws.Range("E44").Value = '  +7.02%  '

$ws.Range("D45").Value = '105.84'
$ws.Range("E45").Value = '  +4.89%  '

$ws.Range("B46").Value = 'TrustWalletToken'
$ws.Range("C46").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D46").Value = '0.8555'
$ws.Range("E46").Value = '  +3.63%  '

$ws.Range("B47").Value = 'PaxDollar'
$ws.Range("C47").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D47").Value = '0.9973'
$ws.Range("E47").Value = '  -0.41%  '

$ws.Range("B48").Value = 'EnergySwap'
$ws.Range("C48").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D48").Value = '9.962'
$ws.Range("E48").Value = '  +3.54%  '

$ws.Range("B49").Value = 'Aptos'
$ws.Range("C49").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D49").Value = '7.516'
$ws.Range("E49").Value = '  +8.66%  '

$ws.Range("B50").Value = 'Decentraland'
$ws.Range("C50").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D50").Value = '0.4277'
$ws.Range("E50").Value = '  +9.47%  '

$ws.Range("B51").Value = 'Elrond'
$ws.Range("C51").Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range("D51").Value = '36.51'
$ws.Range("E51").Value = '  +4.21%  '
